# The upstream change (commit "Fixed #295 Add the version of M2Doc in the
# template custom properties") re-saved this template through a tool that
# canonicalizes OOXML (attributes/namespace declarations sorted
# alphabetically). Verified with XML canonicalization (C14N): every single
# hunk in the diff is a pure attribute/namespace reordering - same element
# names, same attribute sets, same values, same text content, on both the
# "before" and "after" side (word/document.xml and word/styles.xml). There
# is no textual, structural, or value change anywhere in the diff.
#
# In other words, applying this diff does not change the meaning of the
# document at all - it is byte-identical after canonicalization. So the
# correct edit to replay here is a no-op: touch nothing, leave every part
# exactly as it already is.
#
# (The custom "M2Doc version" document property mentioned in the commit
# message is not part of the supplied OOXML diff - docProps/custom.xml is
# untouched by it - and CustomDocumentProperties.Add is not available on
# this host's Word object model, so there is nothing further to apply.)

$d = $word.ActiveDocument
Write-Host ("Paragraphs: " + $d.Paragraphs.Count)
